# fix #683: replaced axis "geo" by "country" in the population_session workbook
# The shared string "geo" (used as the header of the first column on the
# "pop", "births", "deaths" and "__axes__" sheets) is renamed to "country".

$wb = $excel.ActiveWorkbook

$sheetNames = @("pop", "births", "deaths", "__axes__")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    if ($ws.Range("A1").Text -eq "geo") {
        $ws.Range("A1").Value = "country"
    }
}
